$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 3 through 13 (old extra commodity rows), leaving only header + row 2
$ws.Range("A3:B13").EntireRow.Delete()

# Update header text in A1
$ws.Range("A1").Value = "Commodities Down in Price"

# Update row 2 contents: A2 text, B2 numeric value
$ws.Range("A2").Value = "Steel Products (2)."
$ws.Range("B2").Value = 1
